$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that must always be stored as text (strings), even if the text
# looks like a date/number Excel would otherwise auto-convert.
$textCols = @(1,2,3,5,7,9,11)

$data = @(
    @("31.03.2026","10:00","55NM123",5,"11:30",6,"13:45",14,"14:30",1,"Nisa Karaman",9,10),
    @("31.03.2026","10:00","55NM123",5,"11:30",6,"13:45",14,"14:30",1,"Nisa Karaman",9,10),
    @("09.10.1998","10:00","55NM123",5,"11:30",6,"13:45",14,"14:30",1,"Nisa Karaman",9,10),
    @("00.10.1998","10:00","55NM123",5,"11:30",6,"13:45",14,"14:30",1,"Nisa Karaman",9,10),
    @("00.10.42","10:00","55NM123",5,"11:30",6,"13:45",14,"14:30",1,"Nisa Karaman",9,10),
    @("03.05.1979","10:00","55NM123",5,"11:30",6,"13:45",14,"14:30",1,"Nisa Karaman",9,10),
    @("03.05.1979","10:00","55NM123",5,"11:30",6,"13:45",14,"14:30",1,"Nisa Karaman",9,10),
    @("123","2357","2578",853,"2368",4680,"3568",169,"9643",3827,"Sgkhfbll",588,689),
    @("15.02.2025","09:00","68HS574",5,"11:30",6,"13:45",14,"14:30",1,"Melih Karaman",9,10),
    @("11.07.2025","10:00","45HD132",564,"12:00",614,"12:20",617,"14:30",50,"Ela karaman ",116,23)
)

$startRow = 5
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = $j + 1
        $cell = $ws.Cells.Item($row, $col)
        if ($textCols -contains $col) {
            # Force text storage so date/number-looking strings (e.g.
            # "09.10.1998", "03.05.1979", "11.07.2025") aren't silently
            # reinterpreted by Excel as date serial values.
            $cell.NumberFormat = "@"
            $cell.Value = $rowData[$j]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $rowData[$j]
        }
    }
}
